# Updated requirement specifications for the "Web UI" sheet (Jobs section rework).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Web UI")

# --- Make room for the new "Jobs" filter/datagrid block ---
# Everything from the old "ICU" row downwards needs to shift down by 4 rows,
# so insert 4 blank rows right where the old "ICU" row (53) used to be.
$ws.Rows("53:56").Insert()

# --- Clear the old "Ward" related content that is being replaced ---
$ws.Range("C49:H49").ClearContents()
$ws.Range("E50").ClearContents()
$ws.Range("E52").ClearContents()

# --- Rewrite the "Jobs" block with the new filter / datagrid fields ---
$ws.Range("C46").Value = "Filter"
$ws.Range("C47").Value = "Vehicle No"
$ws.Range("C48").Value = "Status"

$ws.Range("C51").Value = "Service Date"
$ws.Range("D51").Value = "Today"
$ws.Range("E51").Value = "All"
$ws.Range("F51").Value = "Select"
$ws.Range("G51").Value = "Similar to graph above the datagrid"

$ws.Range("C53").Value = "Data grid columns"

$ws.Range("C54").Value = "Vehicle No"
$ws.Range("D54").Value = "Service Date"
$ws.Range("E54").Value = "In time"
$ws.Range("F54").Value = "Out Time"
$ws.Range("G54").Value = "Token"
$ws.Range("H54").Value = "Status"
$ws.Range("I54").Value = "Action"

# --- Update the view so it matches where the author was working ---
$ws.Application.ActiveWindow.ScrollRow = 42
$ws.Range("D54").Select()
